# Auto-generated edit script: updates currentAveragePrice / profit figures
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW to match the latest
# scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 835
$ws.Range("I12").Value = 807.5
$ws.Range("K12").Value = 807.5
$ws.Range("M12").Value = -637.5
$ws.Range("H16").Value = 3500
$ws.Range("J16").Value = 3500
$ws.Range("L16").Value = 3500
$ws.Range("N16").Value = -3960
$ws.Range("H62").Value = 3509.4285
$ws.Range("I62").Value = 3378
$ws.Range("K62").Value = 3378
$ws.Range("M62").Value = -2754
$ws.Range("H65").Value = 3509.4285
$ws.Range("I65").Value = 3378
$ws.Range("K65").Value = 16890
$ws.Range("M65").Value = -13770
$ws.Range("H87").Value = 74997.586
$ws.Range("J87").Value = 74997.586
$ws.Range("L87").Value = 74997.586
$ws.Range("N87").Value = -77493.586
$ws.Range("H90").Value = 74997.586
$ws.Range("J90").Value = 74997.586
$ws.Range("L90").Value = 224992.758
$ws.Range("N90").Value = -237472.758
$ws.Range("H112").Value = 2697.1667
$ws.Range("J112").Value = 2622.0625
$ws.Range("L112").Value = 7866.1875
$ws.Range("N112").Value = -10082.1875
$ws.Range("H127").Value = 2502.077
$ws.Range("I127").Value = 2043.9166
$ws.Range("J127").Value = 8000
$ws.Range("K127").Value = 6131.7498
$ws.Range("L127").Value = 24000
$ws.Range("M127").Value = -1171.7498
$ws.Range("N127").Value = -33920
$ws.Range("H132").Value = 3940.7437
$ws.Range("I132").Value = 4167.4546
$ws.Range("K132").Value = 12502.3638
$ws.Range("M132").Value = -9972.363799999999
$ws.Range("H138").Value = 2191.2727
$ws.Range("I138").Value = 955
$ws.Range("J138").Value = 2897.7144
$ws.Range("K138").Value = 2865
$ws.Range("L138").Value = 8693.143199999999
$ws.Range("M138").Value = 2275
$ws.Range("N138").Value = -18973.1432
$ws.Range("H141").Value = 9428.333000000001
$ws.Range("J141").Value = 18042
$ws.Range("L141").Value = 54126
$ws.Range("N141").Value = -64486

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 3000
$ws.Range("I18").Value = 3000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = -2678
$ws.Range("H74").Value = 2980.35
$ws.Range("I74").Value = 3275.5625
$ws.Range("J74").Value = 1799.5
$ws.Range("K74").Value = 3275.5625
$ws.Range("L74").Value = 1799.5
$ws.Range("M74").Value = -2401.5625
$ws.Range("N74").Value = -3547.5
$ws.Range("H77").Value = 2980.35
$ws.Range("I77").Value = 3275.5625
$ws.Range("J77").Value = 1799.5
$ws.Range("K77").Value = 16377.8125
$ws.Range("L77").Value = 8997.5
$ws.Range("M77").Value = -12009.8125
$ws.Range("N77").Value = -17733.5
$ws.Range("H113").Value = 49198.5
$ws.Range("J113").Value = 49198.5
$ws.Range("L113").Value = 49198.5
$ws.Range("N113").Value = -57876.5
$ws.Range("H132").Value = 1697.875
$ws.Range("I132").Value = 1669.1428
$ws.Range("K132").Value = 5007.428400000001
$ws.Range("M132").Value = -2477.428400000001
$ws.Range("H135").Value = 69997
$ws.Range("J135").Value = 69997
$ws.Range("L135").Value = 69997
$ws.Range("N135").Value = -80137

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3146.9
$ws.Range("I134").Value = 3252.2222
$ws.Range("K134").Value = 9756.6666
$ws.Range("M134").Value = -7221.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 9000
$ws.Range("I17").Value = 9000
$ws.Range("K17").Value = 9000
$ws.Range("M17").Value = -8826
$ws.Range("H140").Value = 45379.766
$ws.Range("J140").Value = 45379.766
$ws.Range("L140").Value = 45379.766
$ws.Range("N140").Value = -55739.766

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 9250.083000000001
$ws.Range("I9").Value = 1001
$ws.Range("K9").Value = 3003
$ws.Range("M9").Value = -2779
$ws.Range("H11").Value = 2063.3333
$ws.Range("I11").Value = 2063.3333
$ws.Range("K11").Value = 6189.999899999999
$ws.Range("M11").Value = -6049.999899999999
$ws.Range("H131").Value = 2903.6667
$ws.Range("I131").Value = 2290
$ws.Range("J131").Value = 3026.4
$ws.Range("K131").Value = 6870
$ws.Range("L131").Value = 9079.200000000001
$ws.Range("M131").Value = -1830
$ws.Range("N131").Value = -19159.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 736214.5600000001
$ws.Range("I14").Value = 855583.7
$ws.Range("K14").Value = 855583.7
$ws.Range("M14").Value = -855415.7
$ws.Range("H19").Value = 404404400
$ws.Range("I19").Value = 404404400
$ws.Range("K19").Value = 404404400
$ws.Range("M19").Value = -404404112
$ws.Range("H102").Value = 2749
$ws.Range("I102").Value = 2749
$ws.Range("K102").Value = 2749
$ws.Range("M102").Value = -1127
$ws.Range("H122").Value = 3308.6667
$ws.Range("I122").Value = 1636.8334
$ws.Range("K122").Value = 4910.5002
$ws.Range("M122").Value = -2460.5002
$ws.Range("H132").Value = 3227.9092
$ws.Range("I132").Value = 2890.889
$ws.Range("J132").Value = 4744.5
$ws.Range("K132").Value = 8672.667000000001
$ws.Range("L132").Value = 14233.5
$ws.Range("M132").Value = -6142.667000000001
$ws.Range("N132").Value = -19293.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4999.5
$ws.Range("I7").Value = 4999.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4999.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4887.5
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 2674.625
$ws.Range("I22").Value = 3329.8
$ws.Range("J22").Value = 1582.6666
$ws.Range("K22").Value = 3329.8
$ws.Range("L22").Value = 1582.6666
$ws.Range("M22").Value = -3034.8
$ws.Range("N22").Value = -2172.6666
$ws.Range("H27").Value = 2674.625
$ws.Range("I27").Value = 3329.8
$ws.Range("J27").Value = 1582.6666
$ws.Range("K27").Value = 3329.8
$ws.Range("L27").Value = 1582.6666
$ws.Range("M27").Value = -3222.8
$ws.Range("N27").Value = -1796.6666
$ws.Range("H46").Value = 3625.9092
$ws.Range("I46").Value = 2966.3333
$ws.Range("K46").Value = 2966.3333
$ws.Range("M46").Value = -2778.3333
$ws.Range("H55").Value = 951.26666
$ws.Range("I55").Value = 217
$ws.Range("K55").Value = 217
$ws.Range("M55").Value = -44
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14998.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12528.5
$ws.Range("N126").ClearContents()
